$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy style/format from an existing header cell (F1) to the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the text values (PasteSpecial of formats only shouldn't touch values, but ensure correctness)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# New data cells
$ws.Range("G2").Value = 0.1194315095165318
$ws.Range("H2").Value = 0.9890000000000001
